$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: AD1 = Wins, AE1 = Losses, AF1 = Ties.
# Copy the formatting already used across row 1 (bold, centered, top
# aligned, thin box border) from an existing header cell so the new
# headers reuse the same style as their neighbors instead of a new one.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-51: team record columns Wins / Losses / Ties
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 64  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 98  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
